$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.973.01'
$ws.Range('E2').Value = '  -1.09%  '

$ws.Range('D3').Value = '2.043.20'
$ws.Range('E3').Value = '  -2.37%  '

$ws.Range('E4').Value = '  +0.34%  '

$ws.Range('D5').Value = '''250.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.73%  '

$ws.Range('D6').Value = '''0.666'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.32%  '

$ws.Range('D7').Value = '''58.09'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +7.69%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('E9').Value = '  -2.38%  '

$ws.Range('D10').Value = '''0.384'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.17%  '

$ws.Range('D11').Value = '''0.0785'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.82%  '

$ws.Range('E12').Value = '  +1.94%  '

$ws.Range('D13').Value = '''16.07'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.16%  '

$ws.Range('D14').Value = '2.345.49'
$ws.Range('E14').Value = '  -2.24%  '

$ws.Range('D15').Value = '''0.808'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.98%  '

$ws.Range('D16').Value = '''5.56'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.16%  '

$ws.Range('D17').Value = '2.046.09'
$ws.Range('E17').Value = '  -2.22%  '

$ws.Range('D18').Value = '36.901.81'
$ws.Range('E18').Value = '  -1.12%  '

$ws.Range('D19').Value = '''16.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +14.88%  '

$ws.Range('D20').Value = '''74.63'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.78%  '

$ws.Range('D21').Value = '0.0₃0900'
$ws.Range('E21').Value = '  +5.15%  '

$ws.Range('D22').Value = '''5.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.71%  '

$ws.Range('D23').Value = '''236.48'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.14%  '

$ws.Range('E25').Value = '  -3.63%  '

$ws.Range('D26').Value = '''2.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.21%  '

$ws.Range('D27').Value = '''168.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.48%  '

$ws.Range('D28').Value = '''9.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.45%  '

$ws.Range('D29').Value = '''20.12'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.92%  '

$ws.Range('E30').Value = '  +0.64%  '

$ws.Range('D31').Value = '''1.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.39%  '

$ws.Range('D32').Value = '''4.71'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.31%  '

$ws.Range('D33').Value = '''0.0616'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.99%  '

$ws.Range('D34').Value = '''4.44'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.56%  '

$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').Value = '''1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.12%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '''0.0872'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.84%  '

$ws.Range('D37').Value = '''2.24'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.30%  '

$ws.Range('E38').Value = '  -2.98%  '

$ws.Range('D39').Value = '''0.112'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +16.64%  '

$ws.Range('E40').Value = '  -1.13%  '

$ws.Range('D41').Value = '''17.69'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.09%  '

$ws.Range('D42').Value = '''0.0223'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.63%  '

$ws.Range('D43').Value = '''1.13'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.19%  '

$ws.Range('D44').Value = '''96.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.11%  '

$ws.Range('D45').Value = '''2.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.10%  '

$ws.Range('D46').Value = '''4.61'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +13.48%  '

$ws.Range('D47').Value = '''2.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.63%  '

$ws.Range('D48').Value = '1.281.67'
$ws.Range('E48').Value = '  -4.11%  '

$ws.Range('D49').Value = '''2.89'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.74%  '

$ws.Range('D50').Value = '''6.76'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.19%  '

$ws.Range('D51').Value = '2.231.37'
$ws.Range('E51').Value = '  -2.32%  '
